# Add a new "Estimates" worksheet after the existing Sheet1, populate it
# with the story-point estimate row, and format it (wrap text, column
# widths, row height) to match the target workbook.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the last (currently only) worksheet so
# it lands at the end, becomes sheet index 2 / sheetId 2, and becomes the
# active sheet (mirrors workbookView activeTab="1" / tabSelected moving
# from Sheet1 to the new sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Estimates"

# Content
$newSheet.Range("B2").Value = "Total Story Points Esimates (Including Desing, Cut Effort, DB Design, Testing, Requirement Detailing, Code Review, Bug Fixing, Documentation, Release Notes)"
$newSheet.Range("C2").Value = 314

# Formatting: wrap the label text, set column widths and row height.
$newSheet.Range("B2").WrapText = $true
$newSheet.Columns.Item(2).ColumnWidth = 55.166666666666664
$newSheet.Columns.Item(3).ColumnWidth = 9.709635416666666
$newSheet.Rows.Item(2).RowHeight = 43.5

# Page setup (portrait, matching the target pageSetup orientation).
$newSheet.PageSetup.Orientation = 1

# Selection / activation so the new sheet ends up the visible / active tab
# with C2 selected, matching the authored workbook state.
[void]$newSheet.Range("C2").Select()
[void]$newSheet.Activate()
